# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.091.37'
$ws.Range('E2').Value = '  -0.31%  '
$ws.Range('D3').Value = '1.811.02'
$ws.Range('E3').Value = '  +1.51%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '224.39'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.554'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.54%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '31.90'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.00%  '
$ws.Range('E9').Value = '  +3.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0728'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +11.24%  '
$ws.Range('E11').Value = '  -0.26%  '
$ws.Range('D12').Value = '2.070.14'
$ws.Range('E12').Value = '  +1.48%  '
$ws.Range('D13').Value = '1.806.34'
$ws.Range('E13').Value = '  +0.45%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.82'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.95%  '
$ws.Range('E15').Value = '  +2.14%  '
$ws.Range('D16').Value = '34.042.23'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '4.29'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.11%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.19'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.67%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '249.19'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.95%  '
$ws.Range('D20').Value = '0.0₃0794'
$ws.Range('E20').Value = '  +7.46%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.999'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.01%  '
$ws.Range('B22').Value = 'Avalanche'
$ws.Range('C22').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.95'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.98%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.23'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.24%  '
$ws.Range('E24').Value = '  +0.53%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '159.92'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.79%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '16.56'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.30%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.18'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.74%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0528'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.96%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.75'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.32%  '
$ws.Range('E32').Value = '  +1.51%  '
$ws.Range('E33').Value = '  -0.72%  '
$ws.Range('E34').Value = '  +0.53%  '
$ws.Range('D35').Value = '1.430.04'
$ws.Range('E35').Value = '  -0.60%  '
$ws.Range('E36').Value = '  +0.23%  '
$ws.Range('E37').Value = '  +2.00%  '
$ws.Range('E38').Value = '  +1.14%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.954'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +7.47%  '
$ws.Range('E40').Value = '  -1.69%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '80.76'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.34%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.34'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.16%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.15'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.30%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.02'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.45%  '
$ws.Range('E45').Value = '  +0.51%  '
$ws.Range('E46').Value = '  -2.27%  '
$ws.Range('D47').Value = '1.967.48'
$ws.Range('E47').Value = '  +1.24%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '105.98'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +7.32%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.999'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.13%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '11.77'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.46%  '
$ws.Range('D51').Value = '0.0₆0121'
$ws.Range('E51').Value = '  +2.96%  '
